$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, border, centered) from G1 into H1, then set its text
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data - new "Save" column values (rows 2-11)
$saveValues = @(1, 1, 1, 0, 0, 1, 0, 0, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
